$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5625942945480347
$ws.Range("B1").Value = 0.6419602036476135
$ws.Range("C1").Value = 1.135675549507141
$ws.Range("D1").Value = 1.971238851547241
$ws.Range("E1").Value = 3.833479642868042
